# Recreate baseline and match dac scenario pars
#
# For the CAP, CAP_NEW and INVESTMENT sheets:
#   - the two oldest vintage-year columns (2015, 2020 -> O:P after the
#     shift) are dropped so the sheets line up with REMOVAL (A1:N13,
#     years 2025..2110)
#   - all numeric data cells are reset to 0 (baseline values wiped)
# For the REMOVAL sheet (already A1:N13 / 2025..2110):
#   - all numeric data cells are reset to 0 as well
#
# New year header values shared by every sheet (columns B..N)
$newYears = @(2025,2030,2035,2040,2045,2050,2055,2060,2070,2080,2090,2100,2110)

$wb = $excel.ActiveWorkbook

$sheetsToTrim = @("CAP", "CAP_NEW", "INVESTMENT")
foreach ($sheetName in $sheetsToTrim) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Zero out every data value first (rows 2-13, columns B-P), while the
    # sheet still has its original A1:P13 extent.
    $ws.Range("B2:P13").Value = 0

    # Drop the last two year columns (O, P) so the sheet becomes A1:N13.
    $ws.Range("O1:P13").Delete()

    # Rewrite the header row with the new set of vintage years.
    for ($i = 0; $i -lt $newYears.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item(1, $col).Value = $newYears[$i]
    }
}

# REMOVAL already spans A1:N13 with the correct years - just zero its data.
$wsRemoval = $wb.Worksheets.Item("REMOVAL")
$wsRemoval.Range("B2:N13").Value = 0
